# Updated cryptos list on Mon Jun  3 21:37:08 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row with
# new scraped values. Two rows (45/46) also swap which coin (Arweave /
# Monero) occupies them, including their Link column.
#
# For D-column values that look like plain numbers (e.g. "35.57", "1.00"),
# Excel's COM Value setter auto-coerces the string to a real number (losing
# the original text formatting, e.g. "1.00" -> 1). To keep these as literal
# text - matching how the sheet stores them - we temporarily force the
# cell to Text format ("@") before assigning the value, then restore the
# cell's style to "Normal" so no extra formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.021.81"
$ws.Range("E2").Value = "  +1.76%  "
$ws.Range("D3").Value = "3.763.06"
$ws.Range("E3").Value = "  -0.63%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "624.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.26%  "
$ws.Range("D7").Value = "3.759.76"
$ws.Range("E7").Value = "  -0.68%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("E10").Value = "  +1.36%  "
$ws.Range("E11").Value = "  +3.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.73"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000245"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").Value = "4.403.28"
$ws.Range("E15").Value = "  -0.49%  "
$ws.Range("D16").Value = "3.770.64"
$ws.Range("E16").Value = "  -0.77%  "
$ws.Range("D17").Value = "69.033.69"
$ws.Range("E17").Value = "  +1.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.86%  "
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "467.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.00%  "
$ws.Range("E23").Value = "  +2.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000147"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("E27").Value = "  +3.64%  "
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").Value = "3.916.68"
$ws.Range("E30").Value = "  -0.55%  "
$ws.Range("E31").Value = "  +2.46%  "
$ws.Range("E32").Value = "  +2.52%  "
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.170"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +17.36%  "
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").Value = "3.719.30"
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("E39").Value = "  +2.02%  "
$ws.Range("E40").Value = "  +5.73%  "
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.965"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "153.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.85%  "
$ws.Range("B46").Value = "Arweave"
$ws.Range("C46").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.33%  "
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "46.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.83%  "
$ws.Range("E50").Value = "  +1.36%  "
$ws.Range("E51").Value = "  -0.09%  "
